$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 528, pushing existing rows 528.. down by one
# (this matches dimension growing from A1:R567 to A1:R568 and every row n (n>=529)
# in the new file being identical to old row n-1).
$ws.Rows("528:528").Insert()

# Populate the newly inserted row 528 with its data. Columns A, B, C, E, F, G, H,
# I, R keep the same values as the (now shifted) surrounding rows for this
# market/category block; D, J, K, L, M, N, O, P, Q carry the new record's data.
$ws.Range("A528").Value = 10
$ws.Range("B528").Value = "Vega Modelo de Temuco"
$ws.Range("C528").Value = "La Araucanía"
$ws.Range("D528").Value = 45265
$ws.Range("E528").Value = 9
$ws.Range("F528").Value = 100114013
$ws.Range("G528").Value = "Zanahoria"
$ws.Range("H528").Value = "Sin especificar"
$ws.Range("I528").Value = "Primera"
$ws.Range("J528").Value = 50
$ws.Range("K528").Value = 6000
$ws.Range("L528").Value = 6000
$ws.Range("M528").Value = 6000
$ws.Range("N528").Value = "$/saco 25 kilos"
$ws.Range("O528").Value = "Región de La Araucanía"
$ws.Range("P528").Value = 240
$ws.Range("Q528").Value = 25
$ws.Range("R528").Value = "Hortaliza"

# Match the date cell's number format style (style index 2 / numFmtId 165) used
# throughout column D, same as the rows above/below it.
$ws.Range("D528").NumberFormat = "YYYY-MM-DD HH:MM:SS"
